# Update the "Pet" vendor sheet: the "Catnip" row (row 6) ran out of stock
# and the price per unit was updated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Pet")
$ws.Activate()

$ws.Range("B6").Value = 0      # Stock
$ws.Range("C6").Value = 5.99   # Price

# Move the active selection onto the price cell that was just edited.
$ws.Range("C6").Select()
